$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35.475.48"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.19%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.905.18"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.65%  "
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.637"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.14%  "
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "42.09"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.72%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.339"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.62%  "
$ws.Range("E10").Value = "  +1.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0998"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.03%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.182.08"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.57%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "12.32"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +7.92%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.696"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.28%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.880.12"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.21%  "
$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.83"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.84%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "35.535.17"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.40%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "71.85"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0825"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "243.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.89%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.66"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.85"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.04%  "
$ws.Range("E23").Value = "  +0.28%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.29"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "172.27"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.38%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.21"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +19.69%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.58"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +8.32%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.98"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.70%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.126"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.981"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +25.51%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0571"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.45%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.11"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.74%  "
$ws.Range("E33").Value = "  +0.14%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.17"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.11%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.76"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.38%  "
$ws.Range("E36").Value = "  -0.40%  "
$ws.Range("E37").Value = "  +9.12%  "
$ws.Range("E38").Value = "  +2.47%  "
$ws.Range("E39").Value = "  +1.42%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0637"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +16.63%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "90.97"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.75%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "15.73"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.34%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.350.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("B44").Value = "MultiversX"
$ws.Range("C44").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "50.04"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +44.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.38"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.37%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.84"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.24%  "
$ws.Range("E47").Value = "  +0.25%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.76"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.71%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.63"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.092.37"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.66%  "
$ws.Range("E51").Value = "  +1.56%  "
